# Cap column width at 255
# Adds a new worksheet "Absurdly wide column" at the end of the workbook that
# demonstrates ClosedXML capping an absurdly long auto-fit column width at
# the OOXML maximum of 255 characters.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it lands at the end of
# the sheet tab order (matches "sheetId 7" / the last entry in <x:sheets>).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Absurdly wide column"

# A short string in column A and a very long one (which would normally force
# a column width far beyond 255 characters) in column B.
$newSheet.Range("A1").Value = "Some string"
$newSheet.Range("B1").Value = "Lorem Ipsum is simply dummy text of the printing and typesetting industry. Lorem Ipsum has been the industry's standard dummy text ever since the 1500s, when an unknown printer took a galley of type and scrambled it to make a type specimen book. It has survived not only five centuries, but also the leap into electronic typesetting, remaining essentially unchanged. It was popularised in the 1960s with the release of Letraset sheets containing Lorem Ipsum passages, and more recently with desktop publishing software like Aldus PageMaker including versions of Lorem Ipsum."

# Column A is sized to fit "Some string" (~11.92 characters wide).
$newSheet.Columns.Item(1).ColumnWidth = 11.1665
# Column B would auto-fit far wider than the Lorem Ipsum text warrants, but
# is capped at the maximum allowed column width of 255 characters.
$newSheet.Columns.Item(2).ColumnWidth = 254.1665

# Adding the sheet makes it the active one; restore the original active
# sheet ("Adjust To Contents", the first tab) so the workbook-level view
# state is left untouched, as in the source change.
$wb.Worksheets.Item(1).Activate()
